$wb = $excel.ActiveWorkbook

# The live "Feb" sheet is about to be refreshed with this period's numbers.
# Before doing that, snapshot its current state into a brand-new sheet
# (placed at the end of the workbook) so the prior figures are preserved.
$feb = $wb.Worksheets.Item("Feb")
$feb.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$snapshot = $wb.Worksheets.Item($wb.Worksheets.Count)
$snapshot.Name = "Sheet1"

# Roll the snapshot sheet's values back to what they were before the
# latest "Feb" refresh (the live "Feb" sheet keeps today's numbers).
$snapshot.Range("C5").Value = 2
$snapshot.Range("B11").Value = 20
$snapshot.Range("F11").Value = 50

# Leave the new sheet's selection where it was last left, and make it the
# active/visible tab (mirrors it being the most-recently-touched sheet).
[void]$snapshot.Range("H9").Select()
